# Resort the worksheet tabs so that "总计" (the totals/summary sheet) comes
# before "2022-Q1" (the per-quarter fund detail sheet) — matches commit
# "update data with resort sheetname".
#
# Before: [ "2022-Q1" (active tab), "总计" ]
# After:  [ "总计", "2022-Q1" (active tab) ]
#
# Worksheet handles returned by index (Worksheets.Item(N)) in this runtime
# are resolved live/positionally, so once any operation changes sheet
# order we re-resolve worksheets by NAME to stay safe.

$wb = $excel.ActiveWorkbook

$wsQ1 = $wb.Worksheets.Item("2022-Q1")
$wsTotal = $wb.Worksheets.Item("总计")

# Move "总计" so it sits immediately before "2022-Q1" (i.e. becomes first).
$wsTotal.Move($wsQ1)

# Keep "2022-Q1" as the selected/active tab, matching the original workbook
# (it was the active sheet before the reorder too).
$wsQ1 = $wb.Worksheets.Item("2022-Q1")
$wsQ1.Activate()
